$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column updates: force text storage (NumberFormat "@") so numeric-looking
# price strings like "334.07" are not auto-converted to numbers by Excel,
# then restore the default style so no stray formatting is introduced.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.390.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.098.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5211"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4547"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "55.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +15.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08884"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.177"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.101.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.788"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.998"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001144"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06617"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.277"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.466.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.364"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.344.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.204"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1065"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.644"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.372"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.943"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.788"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.36%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06828"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2307"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6851"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.244"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.312"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6336"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.244"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000343"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +17.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.203"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.62%  "
